$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.84 = 27676.31 pesos`n✅ 27676.31 pesos = 6.82 = 954.38 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update rate values on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 146.19
$wsTasas.Range("O10").Value = 4046
$wsTasas.Range("N12").Value = 4059.9
$wsTasas.Range("O12").Value = 140
